# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404 , *_new -> *_FV2410
# Also freeze the header row and wrap the data range in a real table
# (Table1 / ListObject) with an autofilter, mirroring the export change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename the header cells: "<name>_old" -> "<name>_FV2404"
#                              "<name>_new" -> "<name>_FV2410"
# ---------------------------------------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# Column K ("diff") is untouched.

# ---------------------------------------------------------------------
# 2) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------
# 3) Turn the data range into a real table (ListObject) with headers.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A1:U58")
$table = $ws.ListObjects.Add(1, $dataRange, [Type]::Missing, 1)
$table.Name = "Table1"

"done"
